## PQ_Efficiency_Drills.xlsx edit
## Commit message: "Tried something Diramuid appeared to do"
##
## Adds a new worksheet "BYROWnoLambda" (between TLC4Graphics and Formats)
## that demonstrates BYROW() called with a bare range/table instead of a
## LAMBDA-wrapped argument, plus updates the selection/scroll state on the
## now-deselected TLC4Graphics tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new sheet right after "TLC4Graphics" (i.e. before Formats)
# ---------------------------------------------------------------------
$tlc = $wb.Worksheets.Item("TLC4Graphics")
$ws  = $wb.Worksheets.Add($null, $tlc)
$ws.Name = "BYROWnoLambda"

# ---------------------------------------------------------------------
# 2. Header block (A1:C3) — reuse the same named cell styles the other
#    "memo header" sheets (Formats / Lists) already use.
# ---------------------------------------------------------------------
$ws.Range("A1:A3").Style = "Intro_Hd"
$ws.Range("B1:C3").Style = "Intro_Value"

$ws.Range("A1").Value = "FROM:"
$ws.Range("B1").Value = "Mark Biegert"

$ws.Range("A2").Value = "SUBJECT:"
$ws.Range("B2").Value = "BYROW No Lambda"

$ws.Range("A3").Value = "DATE:"
# Force text so "27-Nov-2023" is stored as a literal string (matching the
# other header cells) rather than being auto-parsed into a date serial.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "27-Nov-2023"
$ws.Range("B3").Style = "Intro_Value"

# Note text describing the experiment.
$ws.Range("E4").Value = "I saw Diramuid use a BYROW and tried it with no lambda. I cannot find where that is supported."

# ---------------------------------------------------------------------
# 3. Sample data table (B7:D17) + BYROW formula spilling into F8:F17
# ---------------------------------------------------------------------
$ws.Range("B7").Value = "A"
$ws.Range("C7").Value = "B"
$ws.Range("D7").Value = "C"

$rows = @(
    @(44,29,60),
    @(80,45,18),
    @(80,75,66),
    @(51,27,81),
    @(57,39,16),
    @(42,98,51),
    @(96,88,91),
    @(96, 6,16),
    @(57,23,25),
    @(52,20,75)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 8 + $i
    $ws.Range("B$r").Value = $rows[$i][0]
    $ws.Range("C$r").Value = $rows[$i][1]
    $ws.Range("D$r").Value = $rows[$i][2]
}

$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("B7:D17"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.TableStyle = "Biegert Table Standard"

# BYROW fed the bare table reference (no LAMBDA wrapper) — the point of
# the experiment.
$ws.Range("F8:F17").FormulaArray = "=BYROW(Table3[],LAMBDA(x,SUM(x)))"

$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 4. View / selection state
#    New sheet becomes the active tab; TLC4Graphics keeps a scrolled,
#    non-active view with its own remembered selection.
# ---------------------------------------------------------------------
$tlc.Activate()
$tlc.Range("E18").Select()
$excel.ActiveWindow.ScrollRow = 13

$ws.Activate()
$ws.Range("J12").Select()
